# Adds a new quarterly vintage column (BB) to the forecast grid, and a new
# trailing row (83) for the newest date point - mirrors the target diff
# which appends column BB (a new forecast vintage) across all existing
# rows 1-82, and appends a brand new row 83 (new date + its BB forecast).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (row 1): BB1 date serial, same style as BA1 ---
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# --- BB values for existing rows 2-82 ---
# Rows 2-70: the forecast has already converged, so the new vintage (BB)
# simply repeats the prior vintage's (BA) value.
# Rows 71-82: the newest vintage revises these estimates.
$values = @{
    2 = -0.2522510312815029
    3 = 0.8253228627186644
    4 = -0.8956080111821194
    5 = 0.6802092152583725
    6 = -0.2702467235064887
    7 = 0.2360788655137185
    8 = 0.4710547017577227
    9 = -0.6739669439274962
    10 = -0.2065125907403171
    11 = 0.1861693474447037
    12 = 0.7530542497296722
    13 = 0.3106180987640244
    14 = 0.5709322285556908
    15 = 1.247003875094862
    16 = -0.776883278166693
    17 = 1.384511819460911
    18 = 0.09417825394891111
    19 = -0.1288380179123294
    20 = 0.1460429756367461
    21 = -0.0194429241770564
    22 = 0.2041957853813301
    23 = 0.1772414135220401
    24 = 0.7756459370471021
    25 = 0.5676350012738425
    26 = -0.7653316860800885
    27 = 0.322227538137227
    28 = -0.0579621925135001
    29 = 0.6475935709367775
    30 = 0.9603367340567104
    31 = 0.3700548251239582
    32 = 0.4443178743943008
    33 = 0.5552956277764309
    34 = 0.3556667645181193
    35 = 0.6924895145077272
    36 = 0.3519868436780342
    37 = 0.489196258618918
    38 = 0.1561519231779869
    39 = 0.5297299217112936
    40 = 0.7664191671286744
    41 = 0.04527341468607915
    42 = 0.06335359735614077
    43 = 0.08598758370690973
    44 = 0.2196343350075409
    45 = 0.03810891122928695
    46 = 0.4286160255108911
    47 = 0.9
    48 = 0.2
    49 = 0.2
    50 = 0
    51 = -2.213339122522456
    52 = -11.49785608241407
    53 = 11.22930999924247
    54 = -2.250986781122748
    55 = -5.40098554941693
    56 = 3.942549781810342
    57 = 6.284666192508709
    58 = -1.347757551663406
    59 = 1.396500038188336
    60 = 0.1885473380929312
    61 = 1.224073604180177
    62 = -1.674179157827311
    63 = -0.6601113848982436
    64 = 0.2658000717656392
    65 = -0.01892863903084674
    66 = 0.4261262404008619
    67 = -0.07433428650158191
    68 = 0.09298591595782568
    69 = 0.2693878213604393
    70 = 0.2408767182737677
    71 = 0.5597354586130052
    72 = 0.1368731201391853
    73 = -0.2551464291630765
    74 = 0.4382843088295769
    75 = 0.273497985379609
    76 = 0.0539306827327628
    77 = 0.2306356627296811
    78 = 0.2402234230226422
    79 = 0.1594459847733843
    80 = 0.1941922693254655
    81 = 0.2122850466304503
    82 = 0.1886006359601506
}

foreach ($row in $values.Keys) {
    $cell = "BB" + $row
    $ws.Range($cell).Value = $values[$row]
}

# --- New trailing row 83: new date in column A (same style as A82), and
#     its corresponding BB forecast value ---
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.1922018276915154
